# Scheduled runner update: refresh cached market-board pricing/profit
# columns (H:N) on the Leve profit sheets. Values below are the new
# snapshot pulled for each affected Leve Item ID row, per sheet.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H80").Value = 756.1667
$ws.Range("I80").Value = 396.2
$ws.Range("J80").Value = 1013.2857
$ws.Range("K80").Value = 1188.6
$ws.Range("L80").Value = 3039.8571
$ws.Range("M80").Value = -190.5999999999999
$ws.Range("N80").Value = -5035.8571

$ws.Range("H83").Value = 756.1667
$ws.Range("I83").Value = 396.2
$ws.Range("J83").Value = 1013.2857
$ws.Range("K83").Value = 3565.8
$ws.Range("L83").Value = 9119.5713
$ws.Range("M83").Value = 1426.2
$ws.Range("N83").Value = -19103.5713

$ws.Range("H86").Value = 3000
$ws.Range("I86").Value = 0
$ws.Range("J86").Value = 3000
$ws.Range("K86").Value = 0
$ws.Range("L86").Value = 3000
$ws.Range("N86").Value = -5246

$ws.Range("H89").Value = 3000
$ws.Range("I89").Value = 0
$ws.Range("J89").Value = 3000
$ws.Range("K89").Value = 0
$ws.Range("L89").Value = 15000
$ws.Range("N89").Value = -26232

$ws.Range("H98").Value = 876.3077
$ws.Range("I98").Value = 798.8333
$ws.Range("J98").Value = 1806
$ws.Range("K98").Value = 798.8333
$ws.Range("L98").Value = 1806
$ws.Range("M98").Value = 699.1667
$ws.Range("N98").Value = -4802

$ws.Range("H115").Value = 490.33334
$ws.Range("I115").Value = 490.33334
$ws.Range("J115").Value = 0
$ws.Range("K115").Value = 1471.00002
$ws.Range("L115").Value = 0
$ws.Range("M115").Value = 95.99998000000005

$ws.Range("H122").Value = 876.3077
$ws.Range("I122").Value = 798.8333
$ws.Range("J122").Value = 1806
$ws.Range("K122").Value = 2396.4999
$ws.Range("L122").Value = 5418
$ws.Range("M122").Value = 53.5001000000002
$ws.Range("N122").Value = -10318

$ws.Range("H132").Value = 15869
$ws.Range("I132").Value = 15340.929
$ws.Range("J132").Value = 18333.334
$ws.Range("K132").Value = 46022.787
$ws.Range("L132").Value = 55000.00199999999
$ws.Range("M132").Value = -43492.787
$ws.Range("N132").Value = -60060.00199999999

$ws.Range("H138").Value = 4748.72
$ws.Range("I138").Value = 1810.6
$ws.Range("J138").Value = 5483.25
$ws.Range("K138").Value = 5431.799999999999
$ws.Range("L138").Value = 16449.75
$ws.Range("M138").Value = -291.7999999999993
$ws.Range("N138").Value = -26729.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6528.75
$ws.Range("I32").Value = 4192.2
$ws.Range("J32").Value = 26000
$ws.Range("K32").Value = 4192.2
$ws.Range("L32").Value = 26000
$ws.Range("M32").Value = -3905.2
$ws.Range("N32").Value = -26574

$ws.Range("H110").Value = 2838.4546
$ws.Range("I110").Value = 2088.5
$ws.Range("J110").Value = 4838.3335
$ws.Range("K110").Value = 2088.5
$ws.Range("L110").Value = 4838.3335
$ws.Range("M110").Value = -43.5
$ws.Range("N110").Value = -8928.333500000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 4960.5
$ws.Range("I99").Value = 7910
$ws.Range("J99").Value = 2011
$ws.Range("K99").Value = 7910
$ws.Range("L99").Value = 2011
$ws.Range("M99").Value = -6412
$ws.Range("N99").Value = -5007

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 74.22221999999999
$ws.Range("I7").Value = 71
$ws.Range("J7").Value = 100
$ws.Range("K7").Value = 71
$ws.Range("L7").Value = 100
$ws.Range("M7").Value = 42
$ws.Range("N7").Value = -326

$ws.Range("H16").Value = 1105.4375
$ws.Range("I16").Value = 999.36365
$ws.Range("J16").Value = 1338.8
$ws.Range("K16").Value = 999.36365
$ws.Range("L16").Value = 1338.8
$ws.Range("M16").Value = -712.36365
$ws.Range("N16").Value = -1912.8

$ws.Range("H32").Value = 1037.8182
$ws.Range("I32").Value = 1037.8182
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 1037.8182
$ws.Range("L32").Value = 0
$ws.Range("M32").Value = -721.8181999999999

$ws.Range("H58").Value = 4244.2856
$ws.Range("I58").Value = 574.5
$ws.Range("J58").Value = 5712.2
$ws.Range("K58").Value = 574.5
$ws.Range("L58").Value = 5712.2
$ws.Range("M58").Value = -371.5
$ws.Range("N58").Value = -6118.2

$ws.Range("H105").Value = 1268.8667
$ws.Range("I105").Value = 1317.9231
$ws.Range("J105").Value = 950
$ws.Range("K105").Value = 1317.9231
$ws.Range("L105").Value = 950
$ws.Range("M105").Value = 429.0769
$ws.Range("N105").Value = -4444

$ws.Range("H113").Value = 1105.4375
$ws.Range("I113").Value = 999.36365
$ws.Range("J113").Value = 1338.8
$ws.Range("K113").Value = 999.36365
$ws.Range("L113").Value = 1338.8
$ws.Range("M113").Value = 1170.63635
$ws.Range("N113").Value = -5678.8

$ws.Range("H132").Value = 4894.3
$ws.Range("I132").Value = 4199.8335
$ws.Range("J132").Value = 5936
$ws.Range("K132").Value = 12599.5005
$ws.Range("L132").Value = 17808
$ws.Range("M132").Value = -10069.5005
$ws.Range("N132").Value = -22868

$ws.Range("H134").Value = 1499.7142
$ws.Range("I134").Value = 1516.3334
$ws.Range("J134").Value = 1400
$ws.Range("K134").Value = 4549.0002
$ws.Range("L134").Value = 4200
$ws.Range("M134").Value = -2014.0002
$ws.Range("N134").Value = -9270

$ws.Range("H136").Value = 4244.2856
$ws.Range("I136").Value = 574.5
$ws.Range("J136").Value = 5712.2
$ws.Range("K136").Value = 1723.5
$ws.Range("L136").Value = 17136.6
$ws.Range("M136").Value = 826.5
$ws.Range("N136").Value = -22236.6

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 763.2308
$ws.Range("I97").Value = 717.5
$ws.Range("J97").Value = 915.6667
$ws.Range("K97").Value = 717.5
$ws.Range("L97").Value = 915.6667
$ws.Range("M97").Value = -221.5
$ws.Range("N97").Value = -1907.6667

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 3633.0715
$ws.Range("I82").Value = 2238.3333
$ws.Range("J82").Value = 4679.125
$ws.Range("K82").Value = 2238.3333
$ws.Range("L82").Value = 4679.125
$ws.Range("M82").Value = -1877.3333
$ws.Range("N82").Value = -5401.125

$ws.Range("H85").Value = 3633.0715
$ws.Range("I85").Value = 2238.3333
$ws.Range("J85").Value = 4679.125
$ws.Range("K85").Value = 2238.3333
$ws.Range("L85").Value = 4679.125
$ws.Range("M85").Value = -990.3332999999998
$ws.Range("N85").Value = -7175.125

$ws.Range("H100").Value = 5177.615
$ws.Range("I100").Value = 1909.8572
$ws.Range("J100").Value = 8990
$ws.Range("K100").Value = 1909.8572
$ws.Range("L100").Value = 8990
$ws.Range("M100").Value = -1368.8572
$ws.Range("N100").Value = -10072

$ws.Range("H139").Value = 0
$ws.Range("I139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("K139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("N139").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 1369.125
$ws.Range("I100").Value = 1453.7142
$ws.Range("J100").Value = 777
$ws.Range("K100").Value = 2907.4284
$ws.Range("L100").Value = 1554
$ws.Range("M100").Value = -2366.4284
$ws.Range("N100").Value = -2636

$ws.Range("H107").Value = 842.7778
$ws.Range("I107").Value = 497.66666
$ws.Range("J107").Value = 1533
$ws.Range("K107").Value = 1492.99998
$ws.Range("L107").Value = 4599
$ws.Range("M107").Value = 427.0000199999999
$ws.Range("N107").Value = -8439

$ws.Range("H126").Value = 4422.875
$ws.Range("I126").Value = 1694.6666
$ws.Range("J126").Value = 6059.8
$ws.Range("K126").Value = 5083.9998
$ws.Range("L126").Value = 18179.4
$ws.Range("M126").Value = -2613.9998
$ws.Range("N126").Value = -23119.4

$ws.Range("H132").Value = 5936.75
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 5936.75
$ws.Range("K132").Value = 0
$ws.Range("L132").Value = 17810.25
$ws.Range("N132").Value = -22870.25
$ws.Range("M132").ClearContents()

$ws.Range("H136").Value = 2173.8333
$ws.Range("I136").Value = 2368.6
$ws.Range("J136").Value = 1200
$ws.Range("K136").Value = 7105.799999999999
$ws.Range("L136").Value = 3600
$ws.Range("M136").Value = -4555.799999999999
$ws.Range("N136").Value = -8700
